$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.009891629219055
$ws.Range("B1").Value = 1.088564276695251
$ws.Range("C1").Value = 5.438004493713379
$ws.Range("D1").Value = 1.611319541931152
$ws.Range("E1").Value = 0.981880247592926
